# Update "想去人数" (F column) values on sheet "展览" and "全部类型"
# per the site regeneration commit (456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 68
$ws1.Range("F4").Value = 2047
$ws1.Range("F5").Value = 341
$ws1.Range("F6").Value = 585
$ws1.Range("F8").Value = 2060
$ws1.Range("F9").Value = 10544
$ws1.Range("F10").Value = 178
$ws1.Range("F12").Value = 280
$ws1.Range("F15").Value = 7439
$ws1.Range("F16").Value = 1114
$ws1.Range("F17").Value = 709
$ws1.Range("F18").Value = 210
$ws1.Range("F19").Value = 63
$ws1.Range("F20").Value = 3310

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 68
$ws4.Range("F4").Value = 2047
$ws4.Range("F5").Value = 341
$ws4.Range("F6").Value = 585
$ws4.Range("F9").Value = 2060
$ws4.Range("F11").Value = 2
$ws4.Range("F12").Value = 10544
$ws4.Range("F13").Value = 178
$ws4.Range("F15").Value = 280
$ws4.Range("F18").Value = 7439
$ws4.Range("F19").Value = 1114
$ws4.Range("F20").Value = 709
$ws4.Range("F21").Value = 210
$ws4.Range("F22").Value = 63
$ws4.Range("F23").Value = 3310
